$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.917021666666667
$ws.Range("H2").Value = 5.751065
$ws.Range("I2").Value = 0.004075802778734984
$ws.Range("J2").Value = 0.004163881452308742
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4126116666666667
$ws.Range("N2").Value = 1.237835
$ws.Range("O2").Value = 0.02619655320417558
$ws.Range("P2").Value = 0.0272642491595578
$ws.Range("Q2").Value = 0.7909855049194444
$ws.Range("R2").Value = 7.118869544274999
$ws.Range("S2").Value = 0.0001067719843428577
$ws.Range("T2").Value = 0.0001135251013866069

$ws.Range("G3").Value = 1.917021666666667
$ws.Range("H3").Value = 5.751065
$ws.Range("I3").Value = 0.004075802778734984
$ws.Range("J3").Value = 0.004163881452308742
$ws.Range("O3").Value = 0.02019824545620215
$ws.Range("P3").Value = 0.02102146768743702
$ws.Range("Q3").Value = 0.6098710489177778
$ws.Range("R3").Value = 5.48883944026
$ws.Range("S3").Value = 0.00008232406495595998
$ws.Range("T3").Value = 0.00008753089940402657

$ws.Range("G4").Value = 1.917021666666667
$ws.Range("H4").Value = 5.751065
$ws.Range("I4").Value = 0.004075802778734984
$ws.Range("J4").Value = 0.004163881452308742
$ws.Range("M4").Value = 5.626650333333333
$ws.Range("N4").Value = 16.879951
$ws.Range("O4").Value = 0.357233827170323
$ws.Range("P4").Value = 0.371793647671238
$ws.Range("Q4").Value = 10.78641059975722
$ws.Range("R4").Value = 97.07769539781498
$ws.Range("S4").Value = 0.001456014625438935
$ws.Range("T4").Value = 0.001548104673624479

$ws.Range("G5").Value = 1.917021666666667
$ws.Range("H5").Value = 5.751065
$ws.Range("I5").Value = 0.004075802778734984
$ws.Range("J5").Value = 0.004163881452308742
$ws.Range("M5").Value = 1.85043
$ws.Range("N5").Value = 3.70086
$ws.Range("O5").Value = 0.1174830763686662
$ws.Range("P5").Value = 0.08151423181978301
$ws.Range("Q5").Value = 3.547314402650001
$ws.Range("R5").Value = 21.2838864159
$ws.Range("S5").Value = 0.0004788378491177438
$ws.Range("T5").Value = 0.0003394155979735896

$ws.Range("G6").Value = 1.917021666666667
$ws.Range("H6").Value = 5.751065
$ws.Range("I6").Value = 0.004075802778734984
$ws.Range("J6").Value = 0.004163881452308742
$ws.Range("M6").Value = 7.542782333333332
$ws.Range("N6").Value = 22.628347
$ws.Range("O6").Value = 0.4788882978006332
$ws.Range("P6").Value = 0.4984064036619842
$ws.Range("Q6").Value = 14.45967715995055
$ws.Range("R6").Value = 130.137094439555
$ws.Range("S6").Value = 0.001951854254879487
$ws.Range("T6").Value = 0.00207530517992004

$ws.Range("I7").Value = 0.9237369463641636
$ws.Range("J7").Value = 0.9436990322117234
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4126116666666667
$ws.Range("N7").Value = 1.237835
$ws.Range("O7").Value = 0.02619655320417558
$ws.Range("P7").Value = 0.0272642491595578
$ws.Range("Q7").Value = 179.2683735201194
$ws.Range("R7").Value = 1613.415361681075
$ws.Range("S7").Value = 0.0241987240620915
$ws.Range("T7").Value = 0.02572924554585399

$ws.Range("I8").Value = 0.9237369463641636
$ws.Range("J8").Value = 0.9436990322117234
$ws.Range("O8").Value = 0.02019824545620215
$ws.Range("P8").Value = 0.02102146768743702
$ws.Range("S8").Value = 0.01865786557962602
$ws.Range("T8").Value = 0.01983793871230433

$ws.Range("I9").Value = 0.9237369463641636
$ws.Range("J9").Value = 0.9436990322117234
$ws.Range("M9").Value = 5.626650333333333
$ws.Range("N9").Value = 16.879951
$ws.Range("O9").Value = 0.357233827170323
$ws.Range("P9").Value = 0.371793647671238
$ws.Range("Q9").Value = 2444.624171128877
$ws.Range("R9").Value = 22001.61754015989
$ws.Range("S9").Value = 0.3299900846482975
$ws.Range("T9").Value = 0.3508613054898137

$ws.Range("I10").Value = 0.9237369463641636
$ws.Range("J10").Value = 0.9436990322117234
$ws.Range("M10").Value = 1.85043
$ws.Range("N10").Value = 3.70086
$ws.Range("O10").Value = 0.1174830763686662
$ws.Range("P10").Value = 0.08151423181978301
$ws.Range("Q10").Value = 803.9607292074501
$ws.Range("R10").Value = 4823.7643752447
$ws.Range("S10").Value = 0.1085234582142595
$ws.Range("T10").Value = 0.07692490167981129

$ws.Range("I11").Value = 0.9237369463641636
$ws.Range("J11").Value = 0.9436990322117234
$ws.Range("M11").Value = 7.542782333333332
$ws.Range("N11").Value = 22.628347
$ws.Range("O11").Value = 0.4788882978006332
$ws.Range("P11").Value = 0.4984064036619842
$ws.Range("Q11").Value = 3277.13060475659
$ws.Range("R11").Value = 29494.17544280931
$ws.Range("S11").Value = 0.4423668138598891
$ws.Range("T11").Value = 0.47034564078394

$ws.Range("G12").Value = 1.719022666666667
$ws.Range("H12").Value = 5.157068
$ws.Range("I12").Value = 0.003654834727920005
$ws.Range("J12").Value = 0.003733816222472697
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.4126116666666667
$ws.Range("N12").Value = 1.237835
$ws.Range("O12").Value = 0.02619655320417558
$ws.Range("P12").Value = 0.0272642491595578
$ws.Range("Q12").Value = 0.709288807531111
$ws.Range("R12").Value = 6.38359926778
$ws.Range("S12").Value = 0.00009574407240242498
$ws.Range("T12").Value = 0.0001017996958054945

$ws.Range("G13").Value = 1.719022666666667
$ws.Range("H13").Value = 5.157068
$ws.Range("I13").Value = 0.003654834727920005
$ws.Range("J13").Value = 0.003733816222472697
$ws.Range("O13").Value = 0.02019824545620215
$ws.Range("P13").Value = 0.02102146768743702
$ws.Range("Q13").Value = 0.5468807030524444
$ws.Range("R13").Value = 4.921926327472
$ws.Range("S13").Value = 0.00007382124893638006
$ws.Range("T13").Value = 0.00007849029707153797

$ws.Range("G14").Value = 1.719022666666667
$ws.Range("H14").Value = 5.157068
$ws.Range("I14").Value = 0.003654834727920005
$ws.Range("J14").Value = 0.003733816222472697
$ws.Range("M14").Value = 5.626650333333333
$ws.Range("N14").Value = 16.879951
$ws.Range("O14").Value = 0.357233827170323
$ws.Range("P14").Value = 0.371793647671238
$ws.Range("Q14").Value = 9.672339460407555
$ws.Range("R14").Value = 87.05105514366799
$ws.Range("S14").Value = 0.001305630597529869
$ws.Range("T14").Value = 0.001388209153087166

$ws.Range("G15").Value = 1.719022666666667
$ws.Range("H15").Value = 5.157068
$ws.Range("I15").Value = 0.003654834727920005
$ws.Range("J15").Value = 0.003733816222472697
$ws.Range("M15").Value = 1.85043
$ws.Range("N15").Value = 3.70086
$ws.Range("O15").Value = 0.1174830763686662
$ws.Range("P15").Value = 0.08151423181978301
$ws.Range("Q15").Value = 3.18093111308
$ws.Range("R15").Value = 19.08558667848
$ws.Range("S15").Value = 0.0004293812274550791
$ws.Range("T15").Value = 0.0003043591611311059

$ws.Range("G16").Value = 1.719022666666667
$ws.Range("H16").Value = 5.157068
$ws.Range("I16").Value = 0.003654834727920005
$ws.Range("J16").Value = 0.003733816222472697
$ws.Range("M16").Value = 7.542782333333332
$ws.Range("N16").Value = 22.628347
$ws.Range("O16").Value = 0.4788882978006332
$ws.Range("P16").Value = 0.4984064036619842
$ws.Range("Q16").Value = 12.96621380073289
$ws.Range("R16").Value = 116.695924206596
$ws.Range("S16").Value = 0.001750257581596251
$ws.Range("T16").Value = 0.001860957915377392

$ws.Range("G17").Value = 29.84747
$ws.Range("H17").Value = 59.69494
$ws.Range("I17").Value = 0.06345906427637789
$ws.Range("J17").Value = 0.04322028241076797
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.4126116666666667
$ws.Range("N17").Value = 1.237835
$ws.Range("O17").Value = 0.02619655320417558
$ws.Range("P17").Value = 0.0272642491595578
$ws.Range("Q17").Value = 12.31541434248333
$ws.Range("R17").Value = 73.8924860549
$ws.Range("S17").Value = 0.001662408753603331
$ws.Range("T17").Value = 0.001178368548393631

$ws.Range("G18").Value = 29.84747
$ws.Range("H18").Value = 59.69494
$ws.Range("I18").Value = 0.06345906427637789
$ws.Range("J18").Value = 0.04322028241076797
$ws.Range("O18").Value = 0.02019824545620215
$ws.Range("P18").Value = 0.02102146768743702
$ws.Range("Q18").Value = 9.495514919293335
$ws.Range("R18").Value = 56.97308951576
$ws.Range("S18").Value = 0.00128176175667519
$ws.Range("T18").Value = 0.0009085537701398615

$ws.Range("G19").Value = 29.84747
$ws.Range("H19").Value = 59.69494
$ws.Range("I19").Value = 0.06345906427637789
$ws.Range("J19").Value = 0.04322028241076797
$ws.Range("M19").Value = 5.626650333333333
$ws.Range("N19").Value = 16.879951
$ws.Range("O19").Value = 0.357233827170323
$ws.Range("P19").Value = 0.371793647671238
$ws.Range("Q19").Value = 167.9412770246567
$ws.Range("R19").Value = 1007.64766214794
$ws.Range("S19").Value = 0.022669724400098
$ws.Range("T19").Value = 0.01606902645088047

$ws.Range("G20").Value = 29.84747
$ws.Range("H20").Value = 59.69494
$ws.Range("I20").Value = 0.06345906427637789
$ws.Range("J20").Value = 0.04322028241076797
$ws.Range("M20").Value = 1.85043
$ws.Range("N20").Value = 3.70086
$ws.Range("O20").Value = 0.1174830763686662
$ws.Range("P20").Value = 0.08151423181978301
$ws.Range("Q20").Value = 55.23065391210001
$ws.Range("R20").Value = 220.9226156484
$ws.Range("S20").Value = 0.007455366094665798
$ws.Range("T20").Value = 0.00352306811974783

$ws.Range("G21").Value = 29.84747
$ws.Range("H21").Value = 59.69494
$ws.Range("I21").Value = 0.06345906427637789
$ws.Range("J21").Value = 0.04322028241076797
$ws.Range("M21").Value = 7.542782333333332
$ws.Range("N21").Value = 22.628347
$ws.Range("O21").Value = 0.4788882978006332
$ws.Range("P21").Value = 0.4984064036619842
$ws.Range("Q21").Value = 225.1329694106967
$ws.Range("R21").Value = 1350.79781646418
$ws.Range("S21").Value = 0.03038980327133558
$ws.Range("T21").Value = 0.02154126552160618

$ws.Range("G22").Value = 2.386211
$ws.Range("H22").Value = 7.158633
$ws.Range("I22").Value = 0.005073351852803602
$ws.Range("J22").Value = 0.005182987702727284
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.4126116666666667
$ws.Range("N22").Value = 1.237835
$ws.Range("O22").Value = 0.02619655320417558
$ws.Range("P22").Value = 0.0272642491595578
$ws.Range("Q22").Value = 0.9845784977283333
$ws.Range("R22").Value = 8.861206479554999
$ws.Range("S22").Value = 0.0001329043317354723
$ws.Range("T22").Value = 0.0001413102681180808

$ws.Range("G23").Value = 2.386211
$ws.Range("H23").Value = 7.158633
$ws.Range("I23").Value = 0.005073351852803602
$ws.Range("J23").Value = 0.005182987702727284
$ws.Range("O23").Value = 0.02019824545620215
$ws.Range("P23").Value = 0.02102146768743702
$ws.Range("Q23").Value = 0.7591364410813333
$ws.Range("R23").Value = 6.832227969732
$ws.Range("S23").Value = 0.0001024728060086051
$ws.Range("T23").Value = 0.000108954008517265

$ws.Range("G24").Value = 2.386211
$ws.Range("H24").Value = 7.158633
$ws.Range("I24").Value = 0.005073351852803602
$ws.Range("J24").Value = 0.005182987702727284
$ws.Range("M24").Value = 5.626650333333333
$ws.Range("N24").Value = 16.879951
$ws.Range("O24").Value = 0.357233827170323
$ws.Range("P24").Value = 0.371793647671238
$ws.Range("Q24").Value = 13.42637491855367
$ws.Range("R24").Value = 120.837374266983
$ws.Range("S24").Value = 0.00181237289895868
$ws.Range("T24").Value = 0.001927001903832147

$ws.Range("G25").Value = 2.386211
$ws.Range("H25").Value = 7.158633
$ws.Range("I25").Value = 0.005073351852803602
$ws.Range("J25").Value = 0.005182987702727284
$ws.Range("M25").Value = 1.85043
$ws.Range("N25").Value = 3.70086
$ws.Range("O25").Value = 0.1174830763686662
$ws.Range("P25").Value = 0.08151423181978301
$ws.Range("Q25").Value = 4.41551642073
$ws.Range("R25").Value = 26.49309852438
$ws.Range("S25").Value = 0.0005960329831680395
$ws.Range("T25").Value = 0.0004224872611191965

$ws.Range("G26").Value = 2.386211
$ws.Range("H26").Value = 7.158633
$ws.Range("I26").Value = 0.005073351852803602
$ws.Range("J26").Value = 0.005182987702727284
$ws.Range("M26").Value = 7.542782333333332
$ws.Range("N26").Value = 22.628347
$ws.Range("O26").Value = 0.4788882978006332
$ws.Range("P26").Value = 0.4984064036619842
$ws.Range("Q26").Value = 17.99867017440566
$ws.Range("R26").Value = 161.988031569651
$ws.Range("S26").Value = 0.002429568832932806
$ws.Range("T26").Value = 0.002583234261140595

